{"js": "// The document had two small wording/typo fixes applied to the\n// \"Respuestas teoricas\" specification paragraphs:\n//  1. \"sea e elemento de S1\" -> \"sea e un elemento de S1\" (missing article).\n//  2. \"...S1[i]+1))> con (0 <= i < n)...\" -> drop the stray \"> \" that was\n//     left after \"S1[i]+1))\" so it reads \"...S1[i]+1)) con (0 <= i < n)...\".\n\nconst body = context.document.body;\n\n// --- Change 1: insert the missing word \"un \" before \"elemento de S1\". ---\nconst missingArticle = body.search(\"elemento de S1, si e >\", { matchCase: true });\nmissingArticle.load(\"items\");\nawait context.sync();\n\nif (missingArticle.items.length > 0) {\n  missingArticle.items[0].insertText(\"un \", \"Before\");\n}\n\n// --- Change 2: remove the stray \"&gt;\" between \"S1[i]+1))\" and \" con (\". ---\nconst strayGt = body.search(\"S1[i]+1))> con (0 <= i < \", { matchCase: true });\nstrayGt.load(\"items\");\nawait context.sync();\n\nif (strayGt.items.length > 0) {\n  strayGt.items[0].insertText(\"S1[i]+1)) con (0 <= i < \", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Two small wording/typo fixes in the \"Respuestas teoricas\" specification\n# paragraphs:\n#  1. \"sea e elemento de S1\" -> \"sea e un elemento de S1\" (missing article).\n#  2. \"...S1[i]+1))> con (0 <= i < n)...\" -> drop the stray \"> \" left after\n#     \"S1[i]+1))\" so it reads \"...S1[i]+1)) con (0 <= i < n)...\".\n\n$d = $word.ActiveDocument\n\n# --- Change 1: insert the missing word \"un \" before \"elemento de S1\". ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"elemento de S1, si e >\"\n$find1.MatchCase = $true\n$find1.MatchWildcards = $false\n$find1.Execute() | Out-Null\nif ($find1.Found) {\n    $ip = $find1.Parent.Duplicate\n    $ip.Collapse(1)   # wdCollapseStart\n    $ip.InsertBefore(\"un \")\n}\n\n# --- Change 2: remove the stray \"&gt;\" between \"S1[i]+1))\" and \" con (\". ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"S1[i]+1))> con (0 <= i < \"\n$find2.Replacement.Text = \"S1[i]+1)) con (0 <= i < \"\n$find2.MatchCase = $true\n$find2.MatchWildcards = $false\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n$d.Save()\n"}
